# Monte Carlo for AI notebook created
# Update the probability distribution values in row 2, row 3, row 4, and row 5
# for columns U through AD on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (U2:AD2)
$ws.Range("U2").Value = 0.47
$ws.Range("V2").Value = 0.51
$ws.Range("W2").Value = 0.02
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0
$ws.Range("AC2").Value = 0
$ws.Range("AD2").Value = 0

# Row 3 (U3:AD3)
$ws.Range("U3").Value = 0.25
$ws.Range("V3").Value = 0.44
$ws.Range("W3").Value = 0.17
$ws.Range("X3").Value = 0.07
$ws.Range("Y3").Value = 0.03
$ws.Range("Z3").Value = 0.01
$ws.Range("AA3").Value = 0
$ws.Range("AB3").Value = 0
$ws.Range("AC3").Value = 0
$ws.Range("AD3").Value = 0

# Row 4 (V4:W4 only)
$ws.Range("V4").Value = 0.07
$ws.Range("W4").Value = 0.01

# Row 5 (U5:AC5)
$ws.Range("U5").Value = 0.1
$ws.Range("V5").Value = 0.18
$ws.Range("W5").Value = 0.21
$ws.Range("X5").Value = 0.22
$ws.Range("Y5").Value = 0.1
$ws.Range("Z5").Value = 0.07
$ws.Range("AA5").Value = 0.05
$ws.Range("AB5").Value = 0.03
$ws.Range("AC5").Value = 0.02
